$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  @('000802','Details: 60bpm_mf_44path_Normal_SplitCycles_TimeBased_Upper_Arm. Script used: ResampledVisualisationCycles_XY.  Dataset used: Session01_Exp_F1_001_SavedCycle_4_Resampled.mat. File Location: Visualisations/Session01_IndividualCyclesTimeThreshold. Date Generated: 29-Jun-2023 16:01:49'),
  @('000803','Details: 60bpm_mf_44path_Normal_SplitCycles_TimeBased_Upper_Arm. Script used: ResampledVisualisationCycles_XY.  Dataset used: Session01_Exp_F1_001_SavedCycle_4_Resampled.mat. File Location: Visualisations/Session01_Resampled_Average_Figures. Date Generated: 29-Jun-2023 16:03:24'),
  @('000804','Details: 60bpm_mf_44path_Normal_SplitCycles_TimeBased_Upper_Arm. Script used: ResampledVisualisationCycles_XY.  Dataset used: Session01_Exp_F1_001_SavedCycle_4_Resampled.mat. File Location: Visualisations/Session01_Resampled_Average_Figures. Date Generated: 29-Jun-2023 16:03:58'),
  @('000805','Details: 60bpm_mf_44path_Normal_SplitCycles_TimeBased_Upper_Arm. Script used: ResampledVisualisationCycles_XY.  Dataset used: Session01_Exp_F1_001_SavedCycle_4_Resampled.mat. File Location: Visualisations/Session01_Resampled_Average_Figures. Date Generated: 29-Jun-2023 16:05:55'),
  @('000806','Details: 60bpm_mf_44path_Normal_SplitCycles_TimeBased_Upper_Arm. Script used: ResampledVisualisationCycles_XY.  Dataset used: Session01_Exp_F1_001_SavedCycle_4_Resampled.mat. File Location: Visualisations/Session01_Resampled_Average_Figures. Date Generated: 29-Jun-2023 16:06:14'),
  @('000807','Details: 60bpm_mf_44path_Normal_SplitCycles_TimeBased_Upper_Arm. Script used: ResampledVisualisationCycles_XY.  Dataset used: Session01_Exp_B1_001_SavedCycle_3_Resampled.mat. File Location: Visualisations/Session01_Resampled_Average_Figures. Date Generated: 29-Jun-2023 16:06:56'),
  @('000808','Details: 60bpm_mf_44path_Normal_SplitCycles_TimeBased_Knees. Script used: ResampledVisualisationCycles_XY.  Dataset used: Session01_Exp_B1_001_SavedCycle_3_Resampled.mat. File Location: Visualisations/Session01_Resampled_Average_Figures. Date Generated: 29-Jun-2023 16:07:13'),
  @('000809','Details: 60bpm_mf_44path_Normal_SplitCycles_TimeBased_Waist. Script used: ResampledVisualisationCycles_XY.  Dataset used: Session01_Exp_C1_001_SavedCycle_5_Resampled.mat. File Location: Visualisations/Session01_Resampled_Average_Figures. Date Generated: 29-Jun-2023 16:07:54'),
  @('000810','Details: 60bpm_mf_44path_Normal_SplitCycles_TimeBased_Waist. Script used: ResampledVisualisationCycles_XY.  Dataset used: Session01_Exp_C1_001_SavedCycle_5_Resampled.mat. File Location: Visualisations/Session01_Resampled_Average_Figures. Date Generated: 29-Jun-2023 16:10:20'),
  @('000811','Details: 60bpm_mf_44path_Normal_SplitCycles_TimeBased_Feet. Script used: ResampledVisualisationCycles_XY.  Dataset used: Session01_Exp_D1_002_SavedCycle_4_Resampled.mat. File Location: Visualisations/Session01_Resampled_Average_Figures. Date Generated: 29-Jun-2023 16:11:07'),
  @('000812','Details: 60bpm_mf_44path_Normal_SplitCycles_TimeBased_Feet. Script used: ResampledVisualisationCycles_XY.  Dataset used: Session01_Exp_D1_002_SavedCycle_4_Resampled.mat. File Location: Visualisations/Session01_Resampled_Average_Figures. Date Generated: 29-Jun-2023 16:11:10'),
  @('000813','Details: 60bpm_mf_44path_Normal_SplitCycles_TimeBased_Feet. Script used: ResampledVisualisationCycles_XY.  Dataset used: Session01_Exp_D1_002_SavedCycle_4_Resampled.mat. File Location: Visualisations/Session01_Resampled_Average_Figures. Date Generated: 29-Jun-2023 16:11:41'),
  @('000814','Details: 60bpm_mf_44path_Normal_SplitCycles_TimeBased_Feet. Script used: ResampledVisualisationCycles_XY.  Dataset used: Session01_Exp_D1_002_SavedCycle_4_Resampled.mat. File Location: Visualisations/Session01_Resampled_Average_Figures. Date Generated: 29-Jun-2023 16:11:47'),
  @('000815','Details: 60bpm_mf_44path_Normal_SplitCycles_TimeBased_Feet. Script used: ResampledVisualisationCycles_XY.  Dataset used: Session01_Exp_D1_002_SavedCycle_4_Resampled.mat. File Location: Visualisations/Session01_Resampled_Average_Figures. Date Generated: 29-Jun-2023 16:12:26'),
  @('000816','Details: 60bpm_mf_44path_Normal_SplitCycles_TimeBased_Feet. Script used: ResampledVisualisationCycles_XY.  Dataset used: Session01_Exp_D1_002_SavedCycle_3_Resampled.mat. File Location: Visualisations/Session01_Resampled_Average_Figures. Date Generated: 29-Jun-2023 16:12:33'),
  @('000817','Details: 60bpm_mf_44path_Normal_SplitCycles_TimeBased_Feet. Script used: ResampledVisualisationCycles_XY.  Dataset used: Session01_Exp_D1_001_SavedCycle_4_Resampled.mat. File Location: Visualisations/Session01_Resampled_Average_Figures. Date Generated: 29-Jun-2023 16:12:41'),
  @('000818','Details: 60bpm_mf_44path_Normal_SplitCycles_TimeBased_Feet. Script used: ResampledVisualisationCycles_XY.  Dataset used: Session01_Exp_D1_002_SavedCycle_4_Resampled.mat. File Location: Visualisations/Session01_Resampled_Average_Figures. Date Generated: 29-Jun-2023 16:12:46'),
  @('000819','Details: 60bpm_mf_44path_Normal_SplitCycles_TimeBased_Feet. Script used: ResampledVisualisationCycles_XY.  Dataset used: Session01_Exp_D1_002_SavedCycle_3_Resampled.mat. File Location: Visualisations/Session01_Resampled_Average_Figures. Date Generated: 29-Jun-2023 16:13:09'),
  @('000820','Details: 60bpm_mf_44path_Normal_SplitCycles_TimeBased_Feet. Script used: ResampledVisualisationCycles_XY.  Dataset used: Session01_Exp_D1_002_SavedCycle_3_Resampled.mat. File Location: Visualisations/Session01_Resampled_Average_Figures. Date Generated: 29-Jun-2023 16:13:25'),
  @('000821','Details: 60bpm_mf_44path_Normal_SplitCycles_TimeBased_Waist. Script used: ResampledVisualisationCycles_XY.  Dataset used: Session01_Exp_C1_001_SavedCycle_5_Resampled.mat. File Location: Visualisations/Session01_Resampled_Average_Figures. Date Generated: 29-Jun-2023 16:13:44'),
  @('000822','Details: 60bpm_mf_44path_Normal_SplitCycles_TimeBased_Knees. Script used: ResampledVisualisationCycles_XY.  Dataset used: Session01_Exp_B1_001_SavedCycle_3_Resampled.mat. File Location: Visualisations/Session01_Resampled_Average_Figures. Date Generated: 29-Jun-2023 16:14:02'),
  @('000823','Details: 60bpm_mf_44path_Normal_SplitCycles_TimeBased_Wrist. Script used: ResampledVisualisationCycles_XY.  Dataset used: Session01_Exp_E1_001_SavedCycle_5_Resampled.mat. File Location: Visualisations/Session01_Resampled_Average_Figures. Date Generated: 29-Jun-2023 16:14:29'),
  @('000824','Details: 60bpm_mf_44path_Normal_SplitCycles_TimeBased_Wrist. Script used: ResampledVisualisationCycles_XY.  Dataset used: Session01_Exp_E1_001_SavedCycle_4_Resampled.mat. File Location: Visualisations/Session01_Resampled_Average_Figures. Date Generated: 29-Jun-2023 16:14:50'),
  @('000825','Details: 60bpm_mf_44path_Normal_SplitCycles_TimeBased_Wrist. Script used: ResampledVisualisationCycles_XY.  Dataset used: Session01_Exp_E1_001_SavedCycle_4_Resampled.mat. File Location: Visualisations/Session01_Resampled_Average_Figures. Date Generated: 29-Jun-2023 16:16:10'),
  @('000826','Details: 60bpm_mf_44path_Normal_SplitCycles_TimeBased_Upper_Arm. Script used: ResampledVisualisationCycles_XY.  Dataset used: Session01_Exp_F1_001_SavedCycle_4_Resampled.mat. File Location: Visualisations/Session01_Resampled_Average_Figures. Date Generated: 29-Jun-2023 16:17:06'),
  @('000827','Details: 60bpm_mf_44path_Normal_SplitCycles_TimeBased_Upper_Arm. Script used: ResampledVisualisationCycles_XY.  Dataset used: Session01_Exp_F1_001_SavedCycle_4_Resampled.mat. File Location: Visualisations/Session01_Resampled_Average_Figures. Date Generated: 29-Jun-2023 16:17:14'),
  @('000828','Details: 60bpm_mf_44path_Normal_SplitCycles_TimeBased_Upper_Arm. Script used: ResampledVisualisationCycles_XY.  Dataset used: Session01_Exp_F1_001_SavedCycle_4_Resampled.mat. File Location: Visualisations/Session01_Resampled_Average_Figures. Date Generated: 29-Jun-2023 16:18:14'),
  @('000829','Details: 60bpm_mf_44path_Normal_SplitCycles_TimeBased_Upper_Arm. Script used: ResampledVisualisationCycles_XY.  Dataset used: Session01_Exp_F1_001_SavedCycle_3_Resampled.mat. File Location: Visualisations/Session01_Resampled_Average_Figures. Date Generated: 29-Jun-2023 16:18:28')
)

$startRow = 804
$helperCell = $ws.Range("ZZ1")
$helperCell.NumberFormat = "@"

for ($i = 0; $i -lt $rows.Count; $i++) {
  $r = $startRow + $i
  $guid = $rows[$i][0]
  $details = $rows[$i][1]

  # Column A must stay text (preserve leading zeros) -- write to a text-formatted
  # helper cell, copy it, and paste-values into the target so no new style/number
  # format gets attached to the destination cell.
  $helperCell.Value = $guid
  $helperCell.Copy()
  $ws.Range("A" + $r).PasteSpecial(-4163)

  $ws.Range("B" + $r).Value = $details
}

$helperCell.Clear()
$excel.CutCopyMode = $false

# Re-fit the GUID/Information columns now that more (shorter) rows exist --
# mirrors Excel's AutoFit-driven bestFit width recorded in the workbook.
$ws.Columns.Item(1).ColumnWidth = 6.333333333333333
$ws.Columns.Item(2).ColumnWidth = 524.0
